# Adding rudimentary implementation of chart commentaries.
# Adds a "Commentary" worksheet (after "Summary ") with a region / commentary
# table, and updates the view/selection state on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Leave the Summary sheet's selection/scroll position at its new resting spot
# before we move focus to the new sheet.
$ws1.Range("A1").Select()
$ws1.Range("B9").Select()

# Create the new "Commentary" sheet right after "Summary ".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Commentary"

# Column widths.
$ws2.Columns.Item(1).ColumnWidth = 25.0283400809717
$ws2.Columns.Item(2).ColumnWidth = 34.3279352226721

# Header row.
$ws2.Range("A1").Value = "Region"
$ws2.Range("B1").Value = "Commentary"
$ws2.Range("A1:B1").Font.Bold = $true

# Data rows: region name + commentary text.
$regions = @("Cape York", "Wet Tropics", "Burdekin", "Mackay Whitsundays", "Fitzroy", "Burnett Mary", "GBR")
$commentaries = @(
    "This is a commentary about the Cape York region.",
    "This is a commentary about the Wet Tropics region.",
    "This is a commentary about the Burdekin region.",
    "This is a commentary about the Mackay Whitsundays region.",
    "This is a commentary about the Fitzroy region.",
    "This is a commentary about the Burnett Mary region.",
    "This is a commentary about the Great Barrier Reef region."
)

for ($i = 0; $i -lt $regions.Count; $i++) {
    $row = $i + 2
    $ws2.Range("A$row").Value = $regions[$i]
    $ws2.Range("B$row").Value = $commentaries[$i]
}

$ws2.Range("A2:B8").RowHeight = 56.7
$ws2.Range("A2:A8").VerticalAlignment = -4160
$ws2.Range("B2:B8").VerticalAlignment = -4160
$ws2.Range("B2:B8").WrapText = $true

# Page setup / print options for the new sheet.
$ps = $ws2.PageSetup
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7
$ps.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ps.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
$ps.PrintGridlines = $false
$ps.PrintHeadings = $false
$ps.CenterHorizontally = $false
$ps.CenterVertically = $false
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.Orientation = 1
$ps.PaperSize = 9

# Activate the new sheet and select its last populated cell, making
# "Commentary" the workbook's active tab.
$ws2.Activate()
$ws2.Range("B8").Select()
